$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "GDP per capita" sheet: refresh the BGDP (B column) figures that other
#    workbooks key off of. The dependent C (GDP per capita) and J (ratio)
#    columns are formulas, so they recalc automatically once B is updated.
# ---------------------------------------------------------------------------
$wsGdp = $wb.Worksheets.Item("GDP per capita")

$wsGdp.Range("B3").Value = 3095729204374.6685
$wsGdp.Range("B4").Value = 3334737305177.4668
$wsGdp.Range("B5").Value = 3578326280207.1709
$wsGdp.Range("B6").Value = 3826013248810.376
$wsGdp.Range("B7").Value = 4078086195995.8813
$wsGdp.Range("B8").Value = 4335826158526.9712
$wsGdp.Range("B9").Value = 4600444659544.1143
$wsGdp.Range("B10").Value = 4872282646816.3506
$wsGdp.Range("B11").Value = 5151896229326.1904
$wsGdp.Range("B12").Value = 5439493947942.0752
$wsGdp.Range("B13").Value = 5735181728184.4834
$wsGdp.Range("B14").Value = 6038976120915.9902
$wsGdp.Range("B15").Value = 6350863885446.5352
$wsGdp.Range("B16").Value = 6670590138492.4668
$wsGdp.Range("B17").Value = 6998184671606.4209
$wsGdp.Range("B18").Value = 7333544869440.4346
$wsGdp.Range("B19").Value = 7676468811471.0947
$wsGdp.Range("B20").Value = 8026771128037.5654
$wsGdp.Range("B21").Value = 8384484920864.9961
$wsGdp.Range("B22").Value = 8749692944266.2607
$wsGdp.Range("B23").Value = 9122398508413.873
$wsGdp.Range("B24").Value = 9502535409857.5332
$wsGdp.Range("B25").Value = 9890044065491.9766
$wsGdp.Range("B26").Value = 10285126395840.611
$wsGdp.Range("B27").Value = 10687971080736.789
$wsGdp.Range("B28").Value = 11098601291388.115
$wsGdp.Range("B29").Value = 11516861449686.387
$wsGdp.Range("B30").Value = 11942572806315.797
$wsGdp.Range("B31").Value = 12375857837659.4
$wsGdp.Range("B32").Value = 12816779436994.98
$wsGdp.Range("B33").Value = 13265191956731.877
$wsGdp.Range("B34").Value = 13720790860998.717

# Reflect where the user was last working on this sheet.
[void]$wsGdp.Activate()
[void]$wsGdp.Range("B35:B46").Select()

# ---------------------------------------------------------------------------
# 2. "Capital Equipment Cost Scalar" sheet: the India:US ratio column L7:L22
#    gets re-entered as one fill (Excel collapses the repeated B/G formula
#    into a shared-formula group), and rows 21-22 shrink to their
#    auto-fit height now that the row no longer needs to wrap as tall.
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capital Equipment Cost Scalar")

$wsCap.Range("L7:L22").Formula = "=B7/G7"

$wsCap.Rows.Item(21).RowHeight = 28.5
$wsCap.Rows.Item(22).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 3. Leave the workbook with the "About" sheet selected/active, matching the
#    saved view state.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
[void]$wsAbout.Activate()
[void]$wsAbout.Range("E11").Select()
